# Medição da Sprint 6 realizada
# Update the non-conformity counts + observações for the Sprint-06 checklists
# (GPR, GQA) and leave a trail of view-state changes (selection on GCO,
# active tab on MED) matching what a reviewer would do while filling in
# the checklist in Excel.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# GPR - Gerencia de Projeto
# ---------------------------------------------------------------------
$gpr = $wb.Worksheets.Item("GPR")

# Linha 5 - "Os riscos estao atualizados..." -> 1 nao-conformidade
$gpr.Range("C5").Value = 1
$gpr.Range("E5").Value = "Riscos não foram atualizados"

# Linha 7 - "O cronograma do projeto..." -> 1 nao-conformidade
$gpr.Range("C7").Value = 1
$gpr.Range("E7").Value = "Cronograma não foi definido"

# Linha 10 - "As licoes aprendidas..." -> 1 nao-conformidade
$gpr.Range("C10").Value = 1
$gpr.Range("E10").Value = "Lições aprendidas não foram registradas"
$gpr.Range("E10").WrapText = $true
$gpr.Range("E10").HorizontalAlignment = -4108
$gpr.Range("E10").VerticalAlignment = -4108

# Linha 13 - "As revisoes previstas..." -> 1 nao-conformidade
$gpr.Range("C13").Value = 1
$gpr.Range("E13").Value = "Houve atraso nas revisões"

# ---------------------------------------------------------------------
# GQA - Gerencia de Qualidade
# ---------------------------------------------------------------------
$gqa = $wb.Worksheets.Item("GQA")

# Linha 7 - "O codigo produzido esta em conformidade com a estrutura..."
$gqa.Range("C7").Value = 736
$gqa.Range("E7").Value = "Existem problemas de padronização "
$gqa.Range("E7").WrapText = $true
$gqa.Range("E7").HorizontalAlignment = -4108
$gqa.Range("E7").VerticalAlignment = -4108

# Linha 9 - "As solicitacoes de correcao de codigo estao sendo executadas..."
$gqa.Range("C9").Value = 1
$gqa.Range("E9").Value = "Não foram executadas correções de código"
$gqa.Range("E9").WrapText = $true
$gqa.Range("E9").HorizontalAlignment = -4108
$gqa.Range("E9").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# View state: GCO selection moved while reviewing, MED left as the
# active sheet/tab when the workbook was saved.
# ---------------------------------------------------------------------
$gco = $wb.Worksheets.Item("GCO")
$gco.Select() | Out-Null
$gco.Range("C8").Select() | Out-Null

$med = $wb.Worksheets.Item("MED")
$med.Select() | Out-Null
